$d = $word.ActiveDocument

# Locate the paragraph that contains "12. dispwflstk: ..." and then the
# first of the two empty paragraphs that follow it (the one right before
# "(7) seismic_location"). That empty paragraph is replaced with two new
# paragraphs describing functions 13 (show_spectrogram) and 14
# (ispectrogram/ispectrogram_1), including a bookmark around "seismogram".
$rng = $d.Content
$rng.Find.Execute("12. dispwflstk:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$targetPara = $rng.Paragraphs(1).Next()
$insertRange = $targetPara.Range

$xmlFrag = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Normal"/><w:spacing w:lineRule="auto" w:line="360"/><w:jc w:val="both"/><w:rPr/></w:pPr><w:r><w:rPr><w:rFonts w:cs="Lohit Devanagari"/><w:b w:val="false"/><w:bCs w:val="false"/><w:color w:val="00000A"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-GB" w:eastAsia="zh-CN" w:bidi="hi-IN"/></w:rPr><w:t xml:space="preserve">13. </w:t></w:r><w:r><w:rPr><w:rFonts w:cs="Lohit Devanagari"/><w:b w:val="false"/><w:bCs w:val="false"/><w:color w:val="00000A"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-GB" w:eastAsia="zh-CN" w:bidi="hi-IN"/></w:rPr><w:t>show_spectrogram: display the spectrogram of seismic data</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Normal"/><w:spacing w:lineRule="auto" w:line="360"/><w:jc w:val="both"/><w:rPr/></w:pPr><w:r><w:rPr><w:rFonts w:cs="Lohit Devanagari"/><w:b w:val="false"/><w:bCs w:val="false"/><w:color w:val="00000A"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-GB" w:eastAsia="zh-CN" w:bidi="hi-IN"/></w:rPr><w:t>14. ispectrogram/</w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsia="Lucida Sans Unicode" w:cs="Lohit Devanagari"/><w:b w:val="false"/><w:bCs w:val="false"/><w:color w:val="00000A"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-GB" w:eastAsia="zh-CN" w:bidi="hi-IN"/></w:rPr><w:t>ispectrogram_</w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsia="Lucida Sans Unicode" w:cs="Lohit Devanagari"/><w:b w:val="false"/><w:bCs w:val="false"/><w:color w:val="00000A"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-GB" w:eastAsia="zh-CN" w:bidi="hi-IN"/></w:rPr><w:t xml:space="preserve">1: display the </w:t></w:r><w:bookmarkStart w:id="100" w:name="__DdeLink__579_706560269"/><w:r><w:rPr><w:rFonts w:eastAsia="Lucida Sans Unicode" w:cs="Lohit Devanagari"/><w:b w:val="false"/><w:bCs w:val="false"/><w:color w:val="00000A"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-GB" w:eastAsia="zh-CN" w:bidi="hi-IN"/></w:rPr><w:t>seismogram</w:t></w:r><w:bookmarkEnd w:id="100"/><w:r><w:rPr><w:rFonts w:eastAsia="Lucida Sans Unicode" w:cs="Lohit Devanagari"/><w:b w:val="false"/><w:bCs w:val="false"/><w:color w:val="00000A"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-GB" w:eastAsia="zh-CN" w:bidi="hi-IN"/></w:rPr><w:t xml:space="preserve"> and spectrogram of seismic data</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$insertRange.InsertXML($xmlFrag)

Write-Output "done"
